$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.578.70'
$ws.Range('E2').Value = '  +3.28%  '
$ws.Range('D3').Value = '4.001.80'
$ws.Range('E3').Value = '  +1.73%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '528.39'
$ws.Range('E5').Value = '  +4.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.68'
$ws.Range('E6').Value = '  +1.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.736'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('E10').Value = '  +1.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000344'
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '44.50'
$ws.Range('E12').Value = '  +3.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.69'
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').Value = '4.641.20'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('D15').Value = '4.025.59'
$ws.Range('E15').Value = '  +2.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.32'
$ws.Range('E16').Value = '  +7.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.25'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  -1.73%  '
$ws.Range('D20').Value = '71.491.32'
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '441.72'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.58'
$ws.Range('E22').Value = '  +4.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '93.85'
$ws.Range('E23').Value = '  +6.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.34'
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.26'
$ws.Range('E25').Value = '  +3.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.10'
$ws.Range('E26').Value = '  +6.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.91'
$ws.Range('E27').Value = '  -1.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.96'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '701.36'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.60'
$ws.Range('E30').Value = '  +2.07%  '
$ws.Range('E31').Value = '  +1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.90'
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.97'
$ws.Range('E33').Value = '  +16.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '67.42'
$ws.Range('E34').Value = '  +0.29%  '
$ws.Range('E35').Value = '  +3.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.445'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '41.05'
$ws.Range('E37').Value = '  +1.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.66'
$ws.Range('E38').Value = '  +20.49%  '
$ws.Range('E39').Value = '  +2.16%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0495'
$ws.Range('E41').Value = '  +1.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.91'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.14'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.56'
$ws.Range('E45').Value = '  +6.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.24'
$ws.Range('E46').Value = '  +9.11%  '
$ws.Range('E47').Value = '  +1.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000285'
$ws.Range('E48').Value = '  +21.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.27'
$ws.Range('E49').Value = '  +6.49%  '
$ws.Range('E51').Value = '  -4.17%  '
